$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet updates ---
$wsSchedule.Range("A3").Value = 46037.27083333334
$wsSchedule.Range("C3").Value = 9.5
$wsSchedule.Range("D3").Value = 35.91
$wsSchedule.Range("E3").Value = 674.410308
$wsSchedule.Range("F3").Value = 18.78057109440267
$wsSchedule.Range("A4").Value = 46037.9375
$wsSchedule.Range("C4").Value = 4.5
$wsSchedule.Range("D4").Value = 17.01
$wsSchedule.Range("E4").Value = 511.2335182499999
$wsSchedule.Range("F4").Value = 30.05488055555555
$wsSchedule.Range("E5").Value = 604.13795325
$wsSchedule.Range("F5").Value = 17.75831726190476

# --- Detailed sheet updates ---
$wsDetailed.Range("E15").Value = "ON"
$wsDetailed.Range("B33").Value = 56.98
$wsDetailed.Range("B34").Value = 47.29031
$wsDetailed.Range("B35").Value = 56.33992
$wsDetailed.Range("C35").Value = "historical"
$wsDetailed.Range("B36").Value = 57.84171
$wsDetailed.Range("C36").Value = "historical"
$wsDetailed.Range("B37").Value = 17.64783
$wsDetailed.Range("B38").Value = 64.01355
$wsDetailed.Range("B39").Value = 67.39884000000001
$wsDetailed.Range("B40").Value = 139.76909
$wsDetailed.Range("B41").Value = 158.99
$wsDetailed.Range("B42").Value = 140.35781
$wsDetailed.Range("B43").Value = 101.25
$wsDetailed.Range("B44").Value = 105.79
$wsDetailed.Range("B46").Value = 80.02
$wsDetailed.Range("E46").Value = "OFF"
$wsDetailed.Range("B47").Value = 68.67858
$wsDetailed.Range("B48").Value = 57.13714
$wsDetailed.Range("B49").Value = 59.11344
$wsDetailed.Range("B50").Value = 59.53998
$wsDetailed.Range("B51").Value = 58.21375
$wsDetailed.Range("B52").Value = 57.06003
$wsDetailed.Range("B54").Value = 50.63915
$wsDetailed.Range("B55").Value = 56.98
$wsDetailed.Range("B56").Value = 49.85284
$wsDetailed.Range("B57").Value = 49.86968
$wsDetailed.Range("B58").Value = 57.06003
$wsDetailed.Range("B59").Value = 63.81284
$wsDetailed.Range("B60").Value = 63.70868
$wsDetailed.Range("B61").Value = 75.41647
$wsDetailed.Range("B62").Value = 65
$wsDetailed.Range("B65").Value = 41.05059
$wsDetailed.Range("B67").Value = 44.06624
$wsDetailed.Range("B68").Value = 36.06
$wsDetailed.Range("B70").Value = 36.06
$wsDetailed.Range("B71").Value = 36.06
$wsDetailed.Range("B72").Value = 36.06028
$wsDetailed.Range("B75").Value = 36.07
$wsDetailed.Range("B77").Value = 36.0601
$wsDetailed.Range("B79").Value = 30.39742
$wsDetailed.Range("B80").Value = 17.74367
$wsDetailed.Range("B81").Value = 17.58028
$wsDetailed.Range("B82").Value = 20.11467
$wsDetailed.Range("B83").Value = 24.72629
$wsDetailed.Range("B84").Value = 8.25426
$wsDetailed.Range("B85").Value = -10.16743
$wsDetailed.Range("B87").Value = -3.03151
$wsDetailed.Range("B88").Value = -3.09383
$wsDetailed.Range("B93").Value = 73.20007

Write-Host "Applied 62 cell updates"
